# Continuacion de logica de reclutamiento
# Applies the "Recruitment" sheet data/formatting update:
#   - Adds a new ID column value (A2 = 1, centered)
#   - Centers the rest of the row 2 data cells (B2:J2, L2)
#   - Converts the CONSENT_TO_KEEP_DATA cell (K2) from a boolean TRUE to the
#     text "yes" (keeps its existing TRUE/FALSE-style number format)
#   - Moves the active selection on the Recruitment sheet from K10 to K6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recruitment")

# New leading ID value for the data row, centered like the rest of row 2.
$idCell = $ws.Cells.Item(2, 1)
$idCell.Value = 1
$idCell.HorizontalAlignment = -4108

# Center-align the existing text/number cells in row 2 (B2:J2).
$ws.Range("B2:J2").HorizontalAlignment = -4108

# CONSENT_TO_KEEP_DATA: was boolean TRUE, now the literal text "yes".
$consentCell = $ws.Cells.Item(2, 11)
$consentCell.Value = "yes"
$consentCell.HorizontalAlignment = -4108

# SHORTLIST_NOTE cell keeps its text but also becomes centered.
$ws.Cells.Item(2, 12).HorizontalAlignment = -4108

# Update the saved selection/active cell on the Recruitment sheet.
$ws.Activate() | Out-Null
$ws.Range("K6").Select() | Out-Null
